# Applies updated "K" (strikeout) values to column G of Sheet1,
# replacing the previous "Strike#" derived values with recalculated ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newGValues = @{
    2 = 1
    3 = 0
    4 = 3
    5 = 0
    6 = 0
    7 = 2
    8 = 1
    9 = 3
    10 = 0
    11 = 2
    13 = 0
    14 = 2
    15 = 1
    16 = 1
    17 = 1
    18 = 1
    19 = 1
    20 = 1
    21 = 1
    22 = 0
    23 = 1
    24 = 0
    25 = 1
    26 = 0
    27 = 1
    28 = 2
    29 = 2
    31 = 1
    32 = 3
    33 = 0
    34 = 1
    35 = 1
    36 = 1
    37 = 0
    38 = 0
    39 = 1
    40 = 0
    41 = 0
    42 = 2
    43 = 0
    44 = 1
    45 = 1
    46 = 2
    47 = 1
    49 = 0
    50 = 0
    51 = 0
    52 = 1
    53 = 0
    54 = 2
    55 = 1
    56 = 2
    57 = 2
    58 = 0
    59 = 0
    60 = 0
    61 = 2
    62 = 1
    63 = 0
    64 = 0
    65 = 2
    66 = 0
    67 = 1
    68 = 0
    69 = 0
    70 = 1
    71 = 2
    72 = 1
    73 = 1
    74 = 3
    75 = 1
    77 = 1
    78 = 1
    79 = 1
}

foreach ($row in $newGValues.Keys) {
    $ws.Cells.Item([int]$row, 7).Value = $newGValues[$row]
}
